# Fill in the previously-blank log rows (19-21) with real Date / Start Time /
# End Time entries - the D column holds a shared ABS(C-B) formula that Excel
# recalculates automatically once B/C get values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: 2017-11-18, 17:30 -> 18:05
$ws.Range("A19").Value = 43057
$ws.Range("B19").Value = 0.72916666666666663
$ws.Range("C19").Value = 0.75347222222222221

# Row 20: 2017-11-18, 19:20 -> 20:15
$ws.Range("A20").Value = 43057
$ws.Range("B20").Value = 0.80555555555555547
$ws.Range("C20").Value = 0.84375

# Row 21: 2017-11-18 22:20 -> 2017-11-19 01:30 (crosses midnight, so the
# Start/End cells carry the full date+time serial rather than a bare
# fraction-of-a-day).
$ws.Range("A21").Value = 43057
$ws.Range("B21").Value = 43057.930555555555
$ws.Range("C21").Value = 43058.0625

# Reflect the author's final selection/scroll state.
$ws.Range("G25").Select()
